$d = $word.ActiveDocument

# 1) Update the "Curso (semestre ideal)" line: EB (8) -> EA (6), EQD (9) -> EQD (7)
$d.Content.Find.Execute("EB (8), EP (10), EQD (9)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EA (6), EP (10), EQD (7)", 2)

# 2) Remove the trailing "Requisitos" heading paragraph and its "LOB1008..." bullet paragraph
$count = $d.Paragraphs.Count
$lastTwo = $d.Range($d.Paragraphs.Item($count - 1).Range.Start, $d.Paragraphs.Item($count).Range.End)
$lastTwo.Delete()
